$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7706114053726196
$ws.Range("B1").Value = 2.694979667663574
$ws.Range("C1").Value = 7.751698017120361
$ws.Range("D1").Value = 2.292196273803711
$ws.Range("E1").Value = 1.505246877670288
